# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) for the rows whose
# annotations changed after re-running the dialog-act tagger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;   DAMSL = "sv"; Dialog = "Statement-opinion" },
    @{ Row = 16;  DAMSL = "sv"; Dialog = "Statement-opinion" },
    @{ Row = 22;  DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 27;  DAMSL = "%";  Dialog = "Uninterpretable" },
    @{ Row = 37;  DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 38;  DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 39;  DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 40;  DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 44;  DAMSL = "%";  Dialog = "Uninterpretable" },
    @{ Row = 50;  DAMSL = "sv"; Dialog = "Statement-opinion" },
    @{ Row = 57;  DAMSL = "sv"; Dialog = "Statement-opinion" },
    @{ Row = 68;  DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 70;  DAMSL = "b";  Dialog = "Acknowledge (Backchannel)" },
    @{ Row = 73;  DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 78;  DAMSL = "%";  Dialog = "Uninterpretable" },
    @{ Row = 81;  DAMSL = "sv"; Dialog = "Statement-opinion" },
    @{ Row = 100; DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 102; DAMSL = "b";  Dialog = "Acknowledge (Backchannel)" },
    @{ Row = 111; DAMSL = "%";  Dialog = "Uninterpretable" },
    @{ Row = 126; DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 128; DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 129; DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 133; DAMSL = "ba"; Dialog = "Appreciation" },
    @{ Row = 151; DAMSL = "ba"; Dialog = "Appreciation" },
    @{ Row = 166; DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 181; DAMSL = "b";  Dialog = "Acknowledge (Backchannel)" },
    @{ Row = 182; DAMSL = "b";  Dialog = "Acknowledge (Backchannel)" },
    @{ Row = 183; DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 187; DAMSL = "aa"; Dialog = "Agree/Accept" },
    @{ Row = 189; DAMSL = "sd"; Dialog = "Statement-non-opinion" },
    @{ Row = 192; DAMSL = "aa"; Dialog = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSL
    $ws.Cells.Item($u.Row, 10).Value = $u.Dialog
}
